$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoices")
$ws.Range("A1").Value = "TEST"
